$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SqlIni IP address used on Linux was wrong (127.0.0.1 loopback) -
# replace it with the real LAN address of the SQL box.
$ws.Range("C2").Value = "192.168.1.113"

# Give the IP cell the same "store as text" formatting already used by the
# ServerID cell next to it (prevents Excel from mangling the dotted value).
$ws.Range("C2").NumberFormat = "@"

# Column C now holds a longer string ("192.168.1.113" vs "127.0.0.1"), so
# widen it to fit the new content.
$ws.Columns.Item(3).ColumnWidth = 14.285714285714286

# Leave the cursor parked on the cell that was just edited.
$ws.Range("C2").Select()
